$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 882, pushing the existing rows 882:923 down to 883:924
$ws.Rows("882:882").Insert()

# The date column is stored as literal text (e.g. "2026/12/29"), not a real
# date serial, so force text format before writing to avoid Excel's
# autoconversion of "2026/02/27" into a date serial number.
$ws.Range("A882").NumberFormat = "@"
$ws.Range("A882").Value = "2026/02/27"
# Restore the default (un-styled) look so the new row matches its siblings.
$ws.Range("A882").Style = "Normal"

$ws.Range("B882").Value = "金"
$ws.Range("C882").Value = 1
$ws.Range("D882").Value = 201
